$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Birds")

# Update Dutch bird translations (Power text column, column E).
# The 7 brand-new translations below are written first (in the order they
# must land in xl/sharedStrings.xml as new unique strings); the remainder
# reuse already-existing Dutch strings, which Excel's shared-string table
# will de-duplicate automatically.
$ws.Range("E123").Value = 'kies 1 andere speler. Voor elke actiesteen op zijn [grassland], leg je 1 [egg] op deze vogel.'
$ws.Range("E82").Value = 'leg 1 [egg] af om 2 [card] te pakken.'
$ws.Range("E37").Value = 'leg deze vogel horizontaal neer, zodat deze 2 [wetland]-velden bedekt. Betaal de lagere eikosten.'
$ws.Range("E71").Value = 'stop een [card] uit je hand onder deze vogel weg. Doe je dat, pak dan 1 [seed] uit de voorraad.'
$ws.Range("E233").Value = 'stop een [card] uit je hand onder deze vogel weg. Doe je dat, pak dan 1 [invertebrate] uit de voorraad.'
$ws.Range("E138").Value = 'stop ten hoogste 3 [card] uit je hand onder deze vogel weg. Trek 1 [card] voor elke kaart die je hebt weggestopt.'
$ws.Range("E129").Value = 'steel 1 [fish] uit de voorraad van een andere speler en bewaar het op deze kaart. Hij pakt 1 [die] uit het vogelhuisje.'
$ws.Range("E35").Value = 'pak 1 [seed] uit het vogelhuisje (indien beschikbaar). Je mag die op deze kaart bewaren.'
$ws.Range("E41").Value = 'stop een [card] uit je hand onder deze vogel weg. Doe je dat, trek dan 1 [card].'
$ws.Range("E46").Value = 'leg 1 [egg] op deze vogel.'
$ws.Range("E48").Value = 'stop een [card] uit je hand onder deze vogel weg. Doe je dat, trek dan 1 [card].'
$ws.Range("E54").Value = 'als een andere speler de actie "Eieren leggen" uitvoert, legt deze vogel 1 [egg] op een andere vogel met een [bowl]-nest.'
$ws.Range("E56").Value = 'trek 2 nieuwe bonuskaarten en houd er 1.'
$ws.Range("E58").Value = 'pak 1 [seed] uit de voorraad en bewaar die op deze kaart.'
$ws.Range("E66").Value = 'als deze vogel aan de voorwaarden van het "Einde ronde"-doel voldoet, telt hij dubbel.'
$ws.Range("E85").Value = 'pak 1 [seed] uit de voorraad en bewaar die op deze kaart.'
$ws.Range("E87").Value = 'pak 1 [invertebrate] uit de voorraad.'
$ws.Range("E90").Value = 'stop een [card] uit je hand onder deze vogel weg. Doe je dat, leg dan ook 1 [egg] op deze vogel.'
$ws.Range("E108").Value = 'stop een [card] uit je hand onder deze vogel weg. Doe je dat, pak dan 1 [seed] uit de voorraad.'
$ws.Range("E111").Value = 'als deze vogel aan de voorwaarden van het "Einde ronde"-doel voldoet, telt hij dubbel.'
$ws.Range("E133").Value = 'trek 2 nieuwe bonuskaarten en houd er 1.'
$ws.Range("E140").Value = 'leg 1 [egg] af om 2 [card] te pakken.'
$ws.Range("E153").Value = 'als deze vogel zich rechts van alle andere vogels in zijn leefomgeving bevindt, verplaats deze dan naar een andere leefomgeving.'
$ws.Range("E159").Value = 'werp alle dobbelstenen die niet in het vogelhuisje liggen. Gooi je ten minste 1 [rodent], pak dan 1 [rodent] en bewaar die op deze kaart.'
$ws.Range("E162").Value = 'leg 1 [egg] op een vogel naar keuze.'
$ws.Range("E170").Value = 'trek 2 nieuwe bonuskaarten en houd er 1.'
$ws.Range("E171").Value = 'leg alle resterende open [card] af en vul de vogelhouder aan. Doe je dat, trek dan 1 van de nieuwe open [card].'
$ws.Range("E192").Value = 'pak 1 [seed] uit het vogelhuisje (indien beschikbaar). Je mag die op deze kaart bewaren.'
$ws.Range("E219").Value = 'leg deze vogel horizontaal neer, zodat deze 2 [forest]-velden bedekt. Betaal de lagere eikosten.'
$ws.Range("E231").Value = 'speel een tweede vogel in je [forest]. Betaal de normale kosten ervan.'
$ws.Range("E255").Value = 'leg 1 [seed] af om 2 [card] van de gedekte stapel onder deze vogel weg te stoppen.'
$ws.Range("E258").Value = 'pak 1 [seed] uit de voorraad en bewaar het op deze kaart. Je mag [seed] dat op deze kaart ligt op elk moment gebruiken.'

# Restore the sheet's last-known selection.
$ws.Range("E22").Select()
